$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155; this pushes the existing row 155
# (and everything below it, down through the old row 226) down by one,
# so the sheet's used range grows from A1:R226 to A1:R227.
$ws.Rows(155).Insert()

# Populate the newly inserted row 155 with the new weekly record.
$ws.Cells.Item(155, 1).Value = 3
$ws.Cells.Item(155, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(155, 3).Value = "Coquimbo"
$ws.Cells.Item(155, 4).Value = 44510
$ws.Cells.Item(155, 5).Value = 5
$ws.Cells.Item(155, 6).Value = 100112009
$ws.Cells.Item(155, 7).Value = "Acelga"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 230
$ws.Cells.Item(155, 11).Value = 2000
$ws.Cells.Item(155, 12).Value = 2300
$ws.Cells.Item(155, 13).Value = 2143
$ws.Cells.Item(155, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(155, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(155, 16).Value = 357
$ws.Cells.Item(155, 17).Value = 6
$ws.Cells.Item(155, 18).Value = "Hortaliza"
